# Prefix each worksheet's protocol name onto the "Step..." / label values
# found in column A (rows 2..N) of that sheet, per the commit:
# "fix: unique command names in XLSX - prefix protocol name to each step"

$wb = $excel.ActiveWorkbook

# Sheets that receive the "<sheetname> " prefix on their column-A labels.
$targetSheets = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol",
    "dickpic",
    "boosters"
)

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    # Row 1 is the header ("Name"); data starts at row 2.
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $val = $cell.Value2
        if ($val -ne $null -and $val -ne "") {
            $prefix = $sheetName + " "
            if (-not $val.ToString().StartsWith($prefix)) {
                $cell.Value2 = $prefix + $val
            }
        }
    }
}
